{"js": "// Update the date title and the practice-table answers to the new values.\nconst body = context.document.body;\n\n// --- 1. Update the title paragraph with the new date ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\nif (titlePara.text.trim() === \"2023-12-08 Friday\") {\n  titlePara.insertText(\"2023-12-09 Saturday\", Word.InsertLocation.replace);\n}\n\n// --- 2. Update each cell of the (first) table with the new answers ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// Row-major list of the new values, matching the non-blank rows of the\n// table in reading order (there are 5 data columns per row).\nconst newValues = [\n  [\"95\u00f74=23, 3\", \"43\u00f75=8, 3\", \"15\u00f75=3, 0\", \"38\u00f75=7, 3\", \"27\u00f79=3, 0\"],\n  [\"92\u00f76=15, 2\", \"69\u00f75=13, 4\", \"91\u00f72=45, 1\", \"47\u00f78=5, 7\", \"15\u00f78=1, 7\"],\n  [\"91\u00f76=15, 1\", \"33\u00f78=4, 1\", \"21\u00f77=3, 0\", \"82\u00f77=11, 5\", \"15\u00f78=1, 7\"],\n  [\"81\u00f77=11, 4\", \"30\u00f72=15, 0\", \"67\u00f73=22, 1\", \"40\u00f75=8, 0\", \"21\u00f79=2, 3\"],\n  [\"21\u00f77=3, 0\", \"98\u00f73=32, 2\", \"54\u00f74=13, 2\", \"78\u00f75=15, 3\", \"33\u00f78=4, 1\"],\n];\n\nconst originalValues = table.values;\nlet newValueIndex = 0;\nfor (let r = 0; r < originalValues.length; r++) {\n  const row = originalValues[r];\n  const isBlankRow = row.every((cell) => cell.trim() === \"\");\n  if (isBlankRow) {\n    continue;\n  }\n  const replacementRow = newValues[newValueIndex];\n  newValueIndex++;\n  for (let c = 0; c < row.length; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = replacementRow[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date title and the practice-table answers to the new values.\n$d = $word.ActiveDocument\n\n# --- 1. Update the title paragraph with the new date ---\n$titleRange = $d.Paragraphs.Item(1).Range\nif ($titleRange.Text.TrimEnd(\"`r\", \"`n\") -eq \"2023-12-08 Friday\") {\n    $titleRange.Text = \"2023-12-09 Saturday\"\n}\n\n# --- 2. Update each cell of the (first) table with the new answers ---\n$table = $d.Tables.Item(1)\n\n# Row-major list of the new values for the 5 non-blank data rows\n# (1-based data-row indices within the table are 1, 5, 9, 13, 17).\n$dataRowIndices = @(1, 5, 9, 13, 17)\n$newValues = @(\n    @(\"95\u00f74=23, 3\", \"43\u00f75=8, 3\", \"15\u00f75=3, 0\", \"38\u00f75=7, 3\", \"27\u00f79=3, 0\"),\n    @(\"92\u00f76=15, 2\", \"69\u00f75=13, 4\", \"91\u00f72=45, 1\", \"47\u00f78=5, 7\", \"15\u00f78=1, 7\"),\n    @(\"91\u00f76=15, 1\", \"33\u00f78=4, 1\", \"21\u00f77=3, 0\", \"82\u00f77=11, 5\", \"15\u00f78=1, 7\"),\n    @(\"81\u00f77=11, 4\", \"30\u00f72=15, 0\", \"67\u00f73=22, 1\", \"40\u00f75=8, 0\", \"21\u00f79=2, 3\"),\n    @(\"21\u00f77=3, 0\", \"98\u00f73=32, 2\", \"54\u00f74=13, 2\", \"78\u00f75=15, 3\", \"33\u00f78=4, 1\")\n)\n\nfor ($i = 0; $i -lt $dataRowIndices.Count; $i++) {\n    $row = $dataRowIndices[$i]\n    $rowValues = $newValues[$i]\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $table.Cell($row, $c)\n        $cell.Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
